# Add two new requirement-category blocks ("SceneControl" and
# "ScreenEffects") to the "EventEngine" worksheet, mirroring the layout
# used by all the other category blocks on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EventEngine")

# ---------------------------------------------------------------------
# 1) Copy the formatting (styles/borders/number formats) of existing
#    blocks onto the new row ranges so the new rows look identical to
#    the rest of the sheet. A 16 row block (1 header + 1 column-header +
#    14 data rows) goes to B150:H165, and an 8 row block (1 header + 1
#    column-header + 6 data rows) goes to B167:H174.
# ---------------------------------------------------------------------
$ws.Range("B35:H50").Copy() | Out-Null
$ws.Range("B150").PasteSpecial(-4122) | Out-Null

$ws.Range("B59:H66").Copy() | Out-Null
$ws.Range("B167").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Fix up the few style cells that differ from the generic templates used
# above so every row matches the very last (bottom-border) row style of
# its block, and the "Run ..." rows (which have no D/E/F/G input) still
# carry the "has value" H-style used throughout the workbook.
$ws.Range("B52:G52").Copy() | Out-Null
$ws.Range("B165").PasteSpecial(-4122) | Out-Null
$ws.Range("C52").Copy() | Out-Null
$ws.Range("C165").PasteSpecial(-4122) | Out-Null

$ws.Range("H38").Copy() | Out-Null
$ws.Range("H155:H164").PasteSpecial(-4122) | Out-Null

$ws.Range("H62").Copy() | Out-Null
$ws.Range("H172:H173").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) SceneControl block (rows 150-165)
#
# The shared-string table records new strings in the order they were
# first typed by the original author: all of column C top-to-bottom for
# the block, then all of column G top-to-bottom. Mirror that order so
# the rebuilt xl/sharedStrings.xml matches exactly.
# ---------------------------------------------------------------------
$sceneLabels = @(
    "Parse battle processing",
    "Parse shop processing",
    "Parse name input processing",
    "Parse open menu",
    "Parse open save",
    "Parse game over",
    "Parse return to title",
    "Run battle processing",
    "Run shop processing",
    "Run name input processing",
    "Run open menu",
    "Run open save",
    "Run game over",
    "Run return to title"
)

$sceneTests = @(
    "Testing.Engine.SceneControl.XmlSceneControlParserTests.ParseBattleProcessingTest()",
    "Testing.Engine.SceneControl.XmlSceneControlParserTests.ParseShopProcessingTest()",
    "Testing.Engine.SceneControl.XmlSceneControlParserTests.ParseNameInputProcessingTest()",
    "Testing.Engine.SceneControl.XmlSceneControlParserTests.ParseOpenMenu()",
    "Testing.Engine.SceneControl.XmlSceneControlParserTests.ParseOpenSave()",
    "Testing.Engine.SceneControl.XmlSceneControlParserTests.ParseGameOver()",
    "Testing.Engine.SceneControl.XmlSceneControlParserTests.ParseReturnToTitle()"
)

$ws.Range("C150").Value = "SceneControl"
for ($i = 0; $i -lt $sceneLabels.Length; $i++) {
    $ws.Cells.Item(152 + $i, 3).Value = $sceneLabels[$i]
}
for ($i = 0; $i -lt $sceneTests.Length; $i++) {
    $ws.Cells.Item(152 + $i, 7).Value = $sceneTests[$i]
}

$ws.Range("B150").Value = "CategoryID"
$ws.Range("B151").Value = "Requirement ID"
$ws.Range("C151").Value = "Label"
$ws.Range("D151").Value = "Standalone"
$ws.Range("E151").Value = "Gamepad"
$ws.Range("F151").Value = "Mobile"
$ws.Range("G151").Value = "Test"
$ws.Range("H151").Value = "Completed?"

for ($i = 0; $i -lt $sceneLabels.Length; $i++) {
    $row = 152 + $i
    if ($row -eq 152) {
        $ws.Cells.Item($row, 2).Value = 1
    } else {
        $ws.Cells.Item($row, 2).Formula = "=B" + ($row - 1) + "+1"
    }
    if ($i -lt 7) {
        $ws.Cells.Item($row, 4).Value = "x"
        $ws.Cells.Item($row, 5).Value = "x"
        $ws.Cells.Item($row, 6).Value = "x"
    }
    $ws.Cells.Item($row, 8).Formula = "=IF(IF(D" + $row + "=" + [char]34 + "x" + [char]34 + ",1, 0)+IF(E" + $row + "=" + [char]34 + "x" + [char]34 + ",1,0)+IF(F" + $row + "=" + [char]34 + "x" + [char]34 + ",1,0)+IF(G" + $row + "=" + [char]34 + [char]34 + ",0,1)=4,TRUE,FALSE)"
}

# ---------------------------------------------------------------------
# 3) ScreenEffects block (rows 167-174)
# ---------------------------------------------------------------------
$screenLabels = @(
    "Parse fade screen",
    "Parse tint screen",
    "Parse flash screen",
    "Run fade screen",
    "Run tint screen",
    "Run flash screen"
)

$screenTests = @(
    "Testing.Engine.ScreenEffects.XmlScreenEffectsParserTests.ParseFadeScreenTest()",
    "Testing.Engine.ScreenEffects.XmlScreenEffectsParserTests.ParseTintScreenTest()",
    "Testing.Engine.ScreenEffects.XmlScreenEffectsParserTests.ParseFlashScreenTest()"
)

$ws.Range("C167").Value = "ScreenEffects"
for ($i = 0; $i -lt $screenLabels.Length; $i++) {
    $ws.Cells.Item(169 + $i, 3).Value = $screenLabels[$i]
}
for ($i = 0; $i -lt $screenTests.Length; $i++) {
    $ws.Cells.Item(169 + $i, 7).Value = $screenTests[$i]
}

$ws.Range("B167").Value = "CategoryID"
$ws.Range("B168").Value = "Requirement ID"
$ws.Range("C168").Value = "Label"
$ws.Range("D168").Value = "Standalone"
$ws.Range("E168").Value = "Gamepad"
$ws.Range("F168").Value = "Mobile"
$ws.Range("G168").Value = "Test"
$ws.Range("H168").Value = "Completed?"

for ($i = 0; $i -lt $screenLabels.Length; $i++) {
    $row = 169 + $i
    if ($row -eq 169) {
        $ws.Cells.Item($row, 2).Value = 1
    } else {
        $ws.Cells.Item($row, 2).Formula = "=B" + ($row - 1) + "+1"
    }
    if ($i -lt 3) {
        $ws.Cells.Item($row, 4).Value = "x"
        $ws.Cells.Item($row, 5).Value = "x"
        $ws.Cells.Item($row, 6).Value = "x"
    }
    $ws.Cells.Item($row, 8).Formula = "=IF(IF(D" + $row + "=" + [char]34 + "x" + [char]34 + ",1, 0)+IF(E" + $row + "=" + [char]34 + "x" + [char]34 + ",1,0)+IF(F" + $row + "=" + [char]34 + "x" + [char]34 + ",1,0)+IF(G" + $row + "=" + [char]34 + [char]34 + ",0,1)=4,TRUE,FALSE)"
}

# ---------------------------------------------------------------------
# 4) Merge the two new section-header rows, matching the rest of sheet
# ---------------------------------------------------------------------
$ws.Range("C150:H150").Merge() | Out-Null
$ws.Range("C167:H167").Merge() | Out-Null

# ---------------------------------------------------------------------
# 5) Leave the view roughly where the author left it (scrolled to the
#    new content, with F160 selected).
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A149").Select() | Out-Null
$ws.Range("F160").Select() | Out-Null

$wb.Save()
